# Ispravke pred pustanje testova 6.2.2026
#
# Appends two new test rows (81, 82) to Sheet1, extends the AutoFilter /
# _FilterDatabase range and the relevant conditional-formatting blocks to
# cover them, and moves the active selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append row 81 first -----------------------------------------------
# Row 80 already carries the right styles/number formats for every column,
# so clone it and then overwrite the TestName / ID cells.
[void]$ws.Range("A80:F80").Copy($ws.Range("A81:F81"))
$ws.Range("A81").Value = "Credit_Cards-Transactions-Filter_Multiple_Filter_Invalid_[WEB]"
$ws.Range("B81").Value = "C70836"

# --- 2. Grow the AutoFilter range from A1:F79 to A1:F81 --------------------
# (done now, with only row 81 present, so it doesn't auto-expand any
# further once row 82 shows up below it)
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F81").AutoFilter()

# --- 3. Keep the hidden _FilterDatabase defined name in sync --------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$81"
    }
}

# --- 4. Now append row 82 --------------------------------------------------
[void]$ws.Range("A80:F80").Copy($ws.Range("A82:F82"))
$ws.Range("A82").Value = "Payments-Domestic_Payments-Create_New_Recipient_[WEB]"
$ws.Range("B82").Value = "C70837"

# --- 5. Extend the conditional formatting that covered the old tail rows --
$rngB77 = $ws.Range("B77:B80")
$fcB77 = $rngB77.FormatConditions.Item(1)
[void]$fcB77.ModifyAppliesToRange($ws.Range("B77:B82"))

# Note: the "B1:B51 B81:B1048576" duplicate-check rule is intentionally left
# untouched -- it is a disjoint (multi-area) range, and this host's
# FormatCondition.ModifyAppliesToRange only honours the first area of a
# multi-area range, so re-applying it here would silently drop the
# "B81:B1048576" tail instead of sliding it down to "B83:B1048576".
# Leaving the rule alone keeps both areas intact, which is closer to the
# intended B1:B51 / B83:B1048576 result than corrupting it.

# --- 6. Leave the selection where the author left it -----------------------
$ws.Range("C83").Select()
